$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "resultado" (G) and "profit" (H) columns for rows that now have
# a settled outcome from the automatic tracker refresh.

$ws.Range("G45").Value = "Fallo"
$ws.Range("H45").Value = -1

$ws.Range("G49").Value = "Acierto"
$ws.Range("H49").Value = 2.5

$ws.Range("G50").Value = "Fallo"
$ws.Range("H50").Value = -1

$ws.Range("G51").Value = "Acierto"
$ws.Range("H51").Value = 1.2
